$d = $word.ActiveDocument

# The claim form shows "Defendant ref: <<respondentExternalReference>>" and,
# further down, an unconditional "2nd Defendant ref: <<respondent2ExternalReference>>"
# paragraph for the (optional) second defendant. Make the matching "1st" label
# for the first defendant's reference line conditional too: only show
# "1st <<es_>>" before "Defendant ref:" when a second defendant actually
# exists (i.e. respondent2ExternalReference isn't blank) - mirroring the
# cs_{!isBlank(respondent2ExternalReference)} guard already used at the end
# of that same line.
$rng = $d.Content
$found = $rng.Find.Execute("Defendant ref: <<", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $rng.InsertBefore("<<cs_{!isBlank(respondent2ExternalReference)}>>1st <<es_>>")
}
